$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Dear [Customer], -> Dear Customer,
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[Customer]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Customer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Intro paragraph: the two isolated "[Product]" runs -> "Product"
#    ("...inspection of your [Product]..." and "...that the [Product] is in
#    excellent..."). Both are self-contained bold/red runs, so a scoped
#    replace-all keeps their formatting intact.
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs(3).Range.Duplicate
$introPara.Find.Execute("[Product]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Product", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "1. The physical structure of the [Product] is strong..." -> "Product"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(7).Range.Duplicate
$p1.Find.Execute("[Product]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Product", 1) | Out-Null

# ---------------------------------------------------------------------------
# 4) "2. The [Product] operates effectively under various conditions."
#    This sentence was plain literal text (no isolated run for "[Product]"),
#    so narrow a Range onto "[Product]", replace its text, re-apply
#    bold+red formatting, and wrap it with the relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(8)
$hit = $p2.Range.Duplicate
$hit.Find.ClearFormatting()
$hit.Find.Execute("[Product]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$hit.Text = "Product"
$hit.Font.Bold = $true
$hit.Font.Color = 255
$d.Bookmarks.Add("_GoBack", $hit) | Out-Null

# ---------------------------------------------------------------------------
# 5) "3. The efficiency of the [Product] meets our expectations." -> "Product"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(9).Range.Duplicate
$p3.Find.Execute("[Product]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Product", 1) | Out-Null

# ---------------------------------------------------------------------------
# 6) Re-adding the "_GoBack" bookmark above (step 4) automatically drops it
#    from its old spot (paragraph 10, which is otherwise left as the same
#    empty paragraph it always was - bookmark names are unique document-wide).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 7) "[YourName]" -> "YourName" (drop the surrounding bracket runs only, so
#    the spell-check proofErr markers around "YourName" stay untouched).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$openBracket = $lastPara.Range.Duplicate
$openBracket.Find.Execute("[", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$openBracket.Delete() | Out-Null

$closeBracket = $lastPara.Range.Duplicate
$closeBracket.Find.Execute("]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$closeBracket.Delete() | Out-Null
